$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column N ("additional_columns") before the existing
# is_complete / comments columns, shifting them right to O / P.
$ws.Columns("N").Insert()
$ws.Range("N1").Value = "additional_columns"

# Rename the lookup helper columns (K10/K11) to their new, clearer names.
$ws.Range("K10").Value = "status_lookup"
$ws.Range("K11").Value = "I_am_lookup"

# Record the additional columns available for the conditional lookup in
# the new N10 cell, with each option on its own line.
$ws.Range("N10").Value = "blue" + [char]10 + "orange"
$ws.Range("N10").Font.Name = "Calibri"
$ws.Range("N10").Font.Size = 11
$ws.Range("N10").WrapText = $true

# Grow row 10 so the wrapped text is fully visible.
$ws.Rows("10").RowHeight = 32

# Leave the selection on the newly added cell.
$ws.Range("N10").Select() | Out-Null
